$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values per cluster analysis re-run (Fe-number vs Frost, no SiO2)
$ws.Range("B2").Value = 161
$ws.Range("B3").Value = 154

# Remove the old 3rd cluster row (row 4) entirely - clear difference between area1 and rest
$ws.Rows(4).Delete()
